# Append the new price row (2025-05-22) to the "Prices" sheet, matching the
# style of the existing data (plain text values, default cell formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

$newRow = 82

# Force text entry so numeric/date-looking values ("35.5", "2025-05-22", ...)
# are stored as literal strings instead of being auto-converted to numbers /
# date serials by the smart-entry logic behind Range.Value.
$ws.Range("A" + $newRow + ":J" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-05-22"
$ws.Range("B" + $newRow).Value = "35.5"
$ws.Range("C" + $newRow).Value = "35.4"
$ws.Range("D" + $newRow).Value = "0.94"
$ws.Range("E" + $newRow).Value = "0.258"
$ws.Range("F" + $newRow).Value = "0.09"
$ws.Range("G" + $newRow).Value = "5,362"
$ws.Range("H" + $newRow).Value = "8,027"
$ws.Range("I" + $newRow).Value = "8,077"
$ws.Range("J" + $newRow).Value = "7.2346"

# Drop the temporary "@" number format so the new row's cells end up with the
# same (default/general) style as every other row in the sheet, while the
# values themselves remain the literal text we just entered.
$ws.Range("A" + $newRow + ":J" + $newRow).ClearFormats()
